$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated M&A counts by year (B2:B21)
$values = @(8, 10, 9, 14, 8, 13, 9, 7, 11, 10, 16, 14, 17, 11, 13, 21, 29, 32, 24, 13)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
